$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet from "1205" to "query" (also updates the defined name reference)
$ws.Name = "query"

# Update the selection to a single active cell B13
$ws.Range("B13").Select()
